$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 20, inheriting formatting from the row above (row 19) -
# this is how the date column (B, style s="2", numFmtId 14) and the order-id
# column (C, style s="1") picked up their number formats/fonts in the source
# file, matching the row just above the new entry.
$ws.Rows(20).Insert(-4121, 0) | Out-Null   # xlShiftDown, xlFormatFromLeftOrAbove

# Fill in the new bonus record (row index 18, 2025-06-14, order 82981004255,
# R$250 bonus credited to technician "Araujo").
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = (Get-Date -Year 2025 -Month 6 -Day 14 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C20").Value = 82981004255
$ws.Range("D20").Value = 250
$ws.Range("E20").Value = "Araujo"

# The longer order-id value no longer fits column C's previous width, so
# Excel grows it to fit the new content (mirrors the existing best-fit
# column B already has).
$ws.Columns("C:C").ColumnWidth = 10.75

# After typing the row, the cursor lands one row below the table.
$ws.Range("B21").Select() | Out-Null
